$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "batsman" column value ends in the same non-breaking space used by the
# existing rows' "James Pattinson" entries (row 2/3 and the sheet name).
$batsman = "James Pattinson" + [char]0x00A0

# New row 4: duplicate of the "Abu Dhabi" match (currently row 3)
$ws.Range("A4:K4").NumberFormat = "@"
$ws.Range("A4").Value = " Abu Dhabi"
$ws.Range("B4").Value = " September 19 2020"
$ws.Range("C4").Value = "Super Kings won by 5 wickets (with 4 balls remaining)"
$ws.Range("D4").Value = "Mumbai Indians"
$ws.Range("E4").Value = "Chennai Super Kings"
$ws.Range("F4").Value = $batsman
$ws.Range("G4").Value = "11"
$ws.Range("H4").Value = "8"
$ws.Range("I4").Value = "2"
$ws.Range("J4").Value = "0"
$ws.Range("K4").Value = "137.50"

# New row 5: duplicate of the "Sharjah" match (currently row 2)
$ws.Range("A5:K5").NumberFormat = "@"
$ws.Range("A5").Value = " Sharjah"
$ws.Range("B5").Value = " November 03 2020"
$ws.Range("C5").Value = "Sunrisers won by 10 wickets (with 17 balls remaining)"
$ws.Range("D5").Value = "Mumbai Indians"
$ws.Range("E5").Value = "Sunrisers Hyderabad"
$ws.Range("F5").Value = $batsman
$ws.Range("G5").Value = "4"
$ws.Range("H5").Value = "5"
$ws.Range("I5").Value = "0"
$ws.Range("J5").Value = "0"
$ws.Range("K5").Value = "80.00"
